$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the underlying values used by the formulas in F4 (=F2-F3) and F5 (=F4/F2)
$ws.Range("F2").Value = 345855
$ws.Range("F3").Value = 272317

# Update the active selection to match the author's final cursor position
$ws.Range("G10").Select()
